# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-11-24 01:14:15
#
# For every data row in the "Recorded By" column (G), the order of the
# last two comma-separated entries is swapped (e.g. "a, System" becomes
# "System, a", and "x, y, System" becomes "x, System, y").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Value2

    if ($text -eq $null) { continue }
    if ($text -eq "") { continue }

    $parts = $text -split ", "
    $n = $parts.Count

    if ($n -ge 2) {
        $tmp = $parts[$n - 1]
        $parts[$n - 1] = $parts[$n - 2]
        $parts[$n - 2] = $tmp

        $joined = [string]::Join(", ", $parts)
        $cell.Value = $joined
    }
}
